$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows for 2008年 and 2009年 (rows 2 and 3), shifting all rows below up.
$ws.Range("A2:A3").EntireRow.Delete()
